# Update NATMI LR-pair TPM figures for Hcrt-Hcrtr1 with freshly recomputed
# values (new TPM input), and relabel column-A cluster references that
# shift from "MuSCs" to the newly-introduced "Inflammatory-Mac" cluster.
# (Row 3/5/7 column D keeps displaying "MuSCs" - it refers to a different
# cluster slot that is unaffected by the relabeling.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.39064
$ws.Range("H2").Value = 1.17192
$ws.Range("I2").Value = 0.5296693860025763
$ws.Range("J2").Value = 0.5296693860025762
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.378587
$ws.Range("N2").Value = 1.135761
$ws.Range("O2").Value = 0.9725321726710616
$ws.Range("P2").Value = 0.9725321726710616
$ws.Range("Q2").Value = 0.14789122568
$ws.Range("R2").Value = 1.33102103112
$ws.Range("S2").Value = 0.5151205187664327
$ws.Range("T2").Value = 0.5151205187664326

# --- Row 3 ---
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.39064
$ws.Range("H3").Value = 1.17192
$ws.Range("I3").Value = 0.5296693860025763
$ws.Range("J3").Value = 0.5296693860025762
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01069266666666667
$ws.Range("N3").Value = 0.032078
$ws.Range("O3").Value = 0.02746782732893832
$ws.Range("P3").Value = 0.02746782732893832
$ws.Range("Q3").Value = 0.004176983306666667
$ws.Range("R3").Value = 0.03759284976
$ws.Range("S3").Value = 0.01454886723614355
$ws.Range("T3").Value = 0.01454886723614354

# --- Row 4 ---
$ws.Range("A4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = 0.1178916666666667
$ws.Range("H4").Value = 0.353675
$ws.Range("I4").Value = 0.1598494949266683
$ws.Range("J4").Value = 0.1598494949266683
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.378587
$ws.Range("N4").Value = 1.135761
$ws.Range("O4").Value = 0.9725321726710616
$ws.Range("P4").Value = 0.9725321726710616
$ws.Range("Q4").Value = 0.04463225240833334
$ws.Range("R4").Value = 0.401690271675
$ws.Range("S4").Value = 0.1554587766014046
$ws.Range("T4").Value = 0.1554587766014046

# --- Row 5 ---
$ws.Range("A5").Value = "Inflammatory-Mac"
$ws.Range("G5").Value = 0.1178916666666667
$ws.Range("H5").Value = 0.353675
$ws.Range("I5").Value = 0.1598494949266683
$ws.Range("J5").Value = 0.1598494949266683
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01069266666666667
$ws.Range("N5").Value = 0.032078
$ws.Range("O5").Value = 0.02746782732893832
$ws.Range("P5").Value = 0.02746782732893832
$ws.Range("Q5").Value = 0.001260576294444445
$ws.Range("R5").Value = 0.01134518665
$ws.Range("S5").Value = 0.004390718325263728
$ws.Range("T5").Value = 0.004390718325263727

# --- Row 6 ---
$ws.Range("G6").Value = 0.228985
$ws.Range("H6").Value = 0.686955
$ws.Range("I6").Value = 0.3104811190707554
$ws.Range("J6").Value = 0.3104811190707554
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.378587
$ws.Range("N6").Value = 1.135761
$ws.Range("O6").Value = 0.9725321726710616
$ws.Range("P6").Value = 0.9725321726710616
$ws.Range("Q6").Value = 0.086690744195
$ws.Range("R6").Value = 0.780216697755
$ws.Range("S6").Value = 0.3019528773032243
$ws.Range("T6").Value = 0.3019528773032243

# --- Row 7 ---
$ws.Range("G7").Value = 0.228985
$ws.Range("H7").Value = 0.686955
$ws.Range("I7").Value = 0.3104811190707554
$ws.Range("J7").Value = 0.3104811190707554
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.01069266666666667
$ws.Range("N7").Value = 0.032078
$ws.Range("O7").Value = 0.02746782732893832
$ws.Range("P7").Value = 0.02746782732893832
$ws.Range("Q7").Value = 0.002448460276666667
$ws.Range("R7").Value = 0.02203614249
$ws.Range("S7").Value = 0.00852824176753105
$ws.Range("T7").Value = 0.00852824176753105

